$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = "iIAHu292"
$ws.Range("B2").Value = 23101230
$ws.Range("C2").Value = "mbkyjti20"
$ws.Range("D2").Value = "f96`$Ek!A"
$ws.Range("F2").Value = "plrFDWTQ"
$ws.Range("G2").Value = "Jjvw"

# Row 3
$ws.Range("A3").Value = "cePwf630"
$ws.Range("B3").Value = 23101229
$ws.Range("C3").Value = "lyhikcp97"
$ws.Range("D3").Value = "D!`$zF8m9"
$ws.Range("F3").Value = "PRmkaiUP"
$ws.Range("G3").Value = "PhvR"
